# Update "想去人数" (interest count) values in column F across the
# 展览, 演出, 本地生活 and 全部类型 sheets, matching output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 11499
$ws1.Range("F4").Value  = 626
$ws1.Range("F6").Value  = 1430
$ws1.Range("F8").Value  = 158
$ws1.Range("F9").Value  = 33
$ws1.Range("F10").Value = 1040
$ws1.Range("F11").Value = 601
$ws1.Range("F12").Value = 713
$ws1.Range("F13").Value = 1209
$ws1.Range("F14").Value = 240
$ws1.Range("F15").Value = 972
$ws1.Range("F21").Value = 274
$ws1.Range("F23").Value = 522
$ws1.Range("F25").Value = 198
$ws1.Range("F27").Value = 344

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 151
$ws2.Range("F9").Value  = 53
$ws2.Range("F10").Value = 464

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value  = 137

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 11499
$ws4.Range("F5").Value  = 626
$ws4.Range("F7").Value  = 137
$ws4.Range("F8").Value  = 1430
$ws4.Range("F11").Value = 158
$ws4.Range("F13").Value = 33
$ws4.Range("F14").Value = 1040
$ws4.Range("F15").Value = 601
$ws4.Range("F16").Value = 713
$ws4.Range("F17").Value = 1209
$ws4.Range("F18").Value = 240
$ws4.Range("F19").Value = 972
$ws4.Range("F25").Value = 151
$ws4.Range("F27").Value = 274
$ws4.Range("F32").Value = 522
$ws4.Range("F34").Value = 198
$ws4.Range("F35").Value = 53
$ws4.Range("F37").Value = 464
$ws4.Range("F39").Value = 344
